# Update "想去人数" (column F) counts across sheets to match the
# gh-pages data refresh generated at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 851
$ws.Range("F5").Value = 1175
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 3830
$ws.Range("F8").Value = 2538
$ws.Range("F10").Value = 2394
$ws.Range("F14").Value = 1624
$ws.Range("F15").Value = 642
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 94
$ws.Range("F18").Value = 304
$ws.Range("F21").Value = 262
$ws.Range("F22").Value = 66
$ws.Range("F23").Value = 436
$ws.Range("F24").Value = 25
$ws.Range("F26").Value = 480
$ws.Range("F27").Value = 669
$ws.Range("F28").Value = 84
$ws.Range("F29").Value = 72
$ws.Range("F30").Value = 367
$ws.Range("F31").Value = 36
$ws.Range("F32").Value = 1608
$ws.Range("F33").Value = 860
$ws.Range("F34").Value = 29
$ws.Range("F35").Value = 4
$ws.Range("F36").Value = 917
$ws.Range("F37").Value = 1948
$ws.Range("F38").Value = 224
$ws.Range("F39").Value = 514
$ws.Range("F41").Value = 8
$ws.Range("F42").Value = 583
$ws.Range("F43").Value = 1231
$ws.Range("F44").Value = 31

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 5
$ws.Range("F4").Value = 62

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 851
$ws.Range("F3").Value = 1175
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 12
$ws.Range("F6").Value = 3830
$ws.Range("F7").Value = 2538
$ws.Range("F8").Value = 2394
$ws.Range("F10").Value = 1624
$ws.Range("F12").Value = 642
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = 94
$ws.Range("F15").Value = 304
$ws.Range("F18").Value = 262
$ws.Range("F19").Value = 66
$ws.Range("F20").Value = 436
$ws.Range("F21").Value = 25
$ws.Range("F23").Value = 480
$ws.Range("F24").Value = 669
$ws.Range("F25").Value = 84
$ws.Range("F26").Value = 62
$ws.Range("F29").Value = 72
$ws.Range("F30").Value = 367
$ws.Range("F31").Value = 36
$ws.Range("F32").Value = 1608
$ws.Range("F33").Value = 860
$ws.Range("F34").Value = 29
$ws.Range("F36").Value = 917
$ws.Range("F37").Value = 1948
$ws.Range("F38").Value = 224
$ws.Range("F42").Value = 514
$ws.Range("F44").Value = 8
$ws.Range("F45").Value = 583
$ws.Range("F46").Value = 1231
$ws.Range("F47").Value = 31
